# micro_reg_history_vol.xlsx -- refresh the regression-output table:
# 1) improved formatting (std-errors widened to match new coefficient precision)
# 2) added aggregate vs. idiosyncratic volatility comparison (inciqr columns)
#
# The sheet is a small results grid where every cell (B1:G17) is a shared
# string (this is how Stata's esttab/estout-style exports write numbers, so
# even "40529" or "0.02" are text, not numeric cells). We reproduce that by
# writing literal text into each changed cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values whose new text contains non-numeric characters (***, (), etc.) --
# these stay text automatically when assigned through Range.Value.
$ws.Range("E2:F2").Value = "-0.17***"
$ws.Range("G2").Value = "-0.16***"

$ws.Range("B4:D4").Value = "-0.11***"
$ws.Range("E4").Value = "-0.25***"
$ws.Range("F4").Value = "-0.26***"
$ws.Range("G4").Value = "-0.24***"
$ws.Range("B5:D5").Value = "(0.02)"
$ws.Range("E5:G5").Value = "(0.03)"

$ws.Range("B6:C6").Value = "-0.16***"
$ws.Range("D6").Value = "-0.15***"
$ws.Range("E6").Value = "-0.33***"
$ws.Range("F6").Value = "-0.34***"
$ws.Range("G6").Value = "-0.33***"
$ws.Range("B7:D7").Value = "(0.02)"
$ws.Range("E7:G7").Value = "(0.03)"

$ws.Range("B8:D8").Value = "-0.13***"
$ws.Range("E8:G8").Value = "-0.31***"
$ws.Range("B9:D9").Value = "(0.02)"
$ws.Range("E9:G9").Value = "(0.03)"

$ws.Range("F10").Value = "0.08***"
$ws.Range("G10").Value = "0.04***"
$ws.Range("D11").Value = "(0.00)"

$ws.Range("B15:D15").Value = "(0.02)"
$ws.Range("E15:G15").Value = "(0.03)"

# ---- Values whose new text is a plain number ("-0.01", "0.02", "41422", ...) --
# Excel auto-converts a bare numeric literal assigned to Range.Value into a
# real numeric cell, which would lose the original text-cell semantics.
# Briefly mark the range as Text ("@") so the literal is stored as a shared
# string, then drop the transient number-format again so the cell keeps its
# original (unstyled) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D10") "-0.01"
Set-TextValue $ws.Range("B14:F14") "0.02"
Set-TextValue $ws.Range("G14") "0.03"
Set-TextValue $ws.Range("B16:D16") "41422"
Set-TextValue $ws.Range("E16:G16") "44421"
